$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D edits to be interpreted as literal text (matches existing
# inlineStr cells), then clear the temporary number format so cell styling
# is left exactly as it was (no style index on these data cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '44.101.14'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").Value = '2.256.52'
$ws.Range("E3").Value = '  +2.45%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '98.51'
$ws.Range("D6").Value = '272.86'
$ws.Range("E6").Value = '  +5.85%  '
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +0.92%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.627'
$ws.Range("E9").Value = '  +4.90%  '
$ws.Range("D10").Value = '48.02'
$ws.Range("E10").Value = '  +7.37%  '
$ws.Range("D11").Value = '0.0943'
$ws.Range("E11").Value = '  +2.54%  '
$ws.Range("D12").Value = '8.22'
$ws.Range("E12").Value = '  +14.09%  '
$ws.Range("D13").Value = '0.105'
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").Value = '15.44'
$ws.Range("E14").Value = '  +7.53%  '
$ws.Range("D15").Value = '2.585.42'
$ws.Range("E15").Value = '  +2.20%  '
$ws.Range("D16").Value = '0.832'
$ws.Range("E16").Value = '  +6.16%  '
$ws.Range("D17").Value = '2.248.46'
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").Value = '44.115.40'
$ws.Range("E19").Value = '  +3.68%  '
$ws.Range("E20").Value = '  +4.94%  '
$ws.Range("D21").Value = '70.86'
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").Value = '2.37'
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").Value = '234.60'
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").Value = '9.76'
$ws.Range("E24").Value = '  +7.57%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = '11.44'
$ws.Range("E26").Value = '  +7.21%  '
$ws.Range("D27").Value = '2.50'
$ws.Range("E27").Value = '  +12.24%  '
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("D29").Value = '40.06'
$ws.Range("E29").Value = '  +2.57%  '
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("D31").Value = '173.56'
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").Value = '0.0914'
$ws.Range("E32").Value = '  +6.63%  '
$ws.Range("D33").Value = '21.18'
$ws.Range("E33").Value = '  +3.80%  '
$ws.Range("D34").Value = '5.67'
$ws.Range("E34").Value = '  +6.38%  '
$ws.Range("E35").Value = '  +1.47%  '
$ws.Range("D36").Value = '0.113'
$ws.Range("E36").Value = '  +0.58%  '
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("E38").Value = '  -2.48%  '
$ws.Range("D39").Value = '3.54'
$ws.Range("E39").Value = '  +24.26%  '
$ws.Range("D40").Value = '0.250'
$ws.Range("E40").Value = '  +25.46%  '
$ws.Range("D41").Value = '2.20'
$ws.Range("E41").Value = '  +5.02%  '
$ws.Range("D42").Value = '12.48'
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("D44").Value = '62.14'
$ws.Range("E44").Value = '  -1.48%  '
$ws.Range("D45").Value = '0.103'
$ws.Range("E45").Value = '  +5.49%  '
$ws.Range("D46").Value = '8.46'
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("D47").Value = '100.55'
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").Value = '1.15'
$ws.Range("E48").Value = '  +3.96%  '
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("D50").Value = '0.429'
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("D51").Value = '2.469.03'
$ws.Range("E51").Value = '  +2.10%  '

# Restore original (unstyled) formatting on column D now that the text
# values are safely stored as strings.
$ws.Range("D2:D51").ClearFormats()
